# GPS_ERT.xlsx — "make 2d topo strings only"
#
# Summary of the edit being reproduced:
#  - ERT_A_1m  (sheet1.xml):  fix electrode-2 X/Y coordinates; becomes the active sheet/tab.
#  - ERT_A_5m  (sheet2.xml):  remove 5 extraneous GPS rows (electrodes 2,13,20,22,26).
#  - ERT_C_1m  (sheet4.xml):  fix a fat-fingered X coordinate (extra leading "1").
#  - ERT_C_5m  (sheet5.xml):  fix two fat-fingered Y coordinates; no longer the active tab.
#  - ERT_F_5m  (sheet12.xml): add per-segment distance (E) / distance-per-electrode (F) columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ERT_A_5m: delete the rows for electrodes 2, 13, 20, 22, 26 (old rows 3-7),
# which shifts electrodes 28, 32, 37, 48 up to rows 3-6.
# ---------------------------------------------------------------------------
$wsA5 = $wb.Worksheets.Item("ERT_A_5m")
$wsA5.Rows("3:7").Delete()
$wsA5.Range("M11").Select()

# ---------------------------------------------------------------------------
# ERT_C_1m: correct the typo'd X coordinate on row 3 (extra leading "1").
# ---------------------------------------------------------------------------
$wsC1 = $wb.Worksheets.Item("ERT_C_1m")
$wsC1.Range("B3").Value = 613148.6
$wsC1.Range("D6").Select()

# ---------------------------------------------------------------------------
# ERT_C_5m: correct two typo'd Y coordinates (extra leading digits).
# ---------------------------------------------------------------------------
$wsC5 = $wb.Worksheets.Item("ERT_C_5m")
$wsC5.Range("C5").Value = 6651267
$wsC5.Range("C10").Value = 6651187.2000000002
$wsC5.Range("C5").Select()

# ---------------------------------------------------------------------------
# ERT_F_5m: add columns E (segment distance) and F (distance per electrode
# spacing) for rows 3-8.
# ---------------------------------------------------------------------------
$wsF5 = $wb.Worksheets.Item("ERT_F_5m")
for ($r = 3; $r -le 8; $r++) {
  $prev = $r - 1
  $wsF5.Range("E$r").Formula = "=SQRT((B$r-B$prev)^2+(C$r-C$prev)^2)"
}
for ($r = 3; $r -le 8; $r++) {
  $prev = $r - 1
  $wsF5.Range("F$r").Formula = "=E$r/(A$r-A$prev)"
}
$wsF5.Range("I6").Select()

# ---------------------------------------------------------------------------
# ERT_A_1m: correct electrode-2 X/Y coordinates, then make it the active
# sheet/selection (Excel clears the previous sheet's tabSelected flag and
# resets the workbook's stored activeTab automatically).
# ---------------------------------------------------------------------------
$wsA1 = $wb.Worksheets.Item("ERT_A_1m")
$wsA1.Range("B2").Value = 613087
$wsA1.Range("C2").Value = 6651252
$wsA1.Activate()
$wsA1.Range("H21").Select()

Write-Host "edit applied"
